$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 720 (shifts the old row 720 and everything below it down by one).
$ws.Rows.Item(720).Insert()

# Fill in the new row's data. Column A holds a date formatted as text (e.g. "2026/01/25"),
# so force a Text number format before assigning it, otherwise Excel will auto-convert the
# look-alike date string into a date serial number. Then clear the formatting so the cell
# doesn't retain an extra/explicit style (matching the plain, unstyled data rows around it).
$ws.Range("A720").NumberFormat = "@"
$ws.Range("A720").Value = "2026/01/25"
$ws.Range("A720").ClearFormats()

$ws.Range("B720").Value = "日"
$ws.Range("C720").Value = 8
$ws.Range("D720").Value = 172
